# TC01_Canine_Filter_StageOfDisease-2.xlsx
# Update the Neo4j "file" query text on the startup sheet (row for FilesTab,
# cell B4) to drop the `File Type` and `Breed` columns from the RETURN
# clause, matching the new ICDC Jenkins automation scripts. Then reflect
# the user's resulting selection/scroll position on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newFileQuery = @'

MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE diag.stage_of_disease IN ['II']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newFileQuery

# Reflect the saved selection/view state: active cell B4 with the sheet
# scrolled so row 4 is at the top.
$ws.Range("B4").Select() | Out-Null
